$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 79: add VL(2-0-0) particulars, 2 days, remarks 11/18,25 ---
$ws.Range("B79").Value = "VL(2-0-0)"
$ws.Range("D79").Value = 2

# --- Row 81: year header 2023 (text, quote-prefixed, bold like other year rows) ---
$ws.Range("A81").Value = "'2023"
$ws.Range("A81").Font.Bold = $true

# --- Row 79: remarks (set after 2023 label so shared-string order matches) ---
$ws.Range("K79").Value = "11/18,25"

# --- Row 82: Jan 2023 SL(2-0-0) leave entry ---
$ws.Range("A82").Value = 44927
$ws.Range("B82").Value = "SL(2-0-0)"
$ws.Range("C82").Value = 1.25
$ws.Range("H82").Value = 2
$ws.Range("K82").Value = "1/26,27/2023"

# --- Row 83: Feb 2023 ---
$ws.Range("A83").Value = 44958
$ws.Range("C83").Value = 1.25

# --- Rows 84-127: fill PERIOD date series (1st of month) ---
$ws.Range("A84").Value = 44986
$ws.Range("A85").Value = 45017
$ws.Range("A86").Value = 45047
$ws.Range("A87").Value = 45078
$ws.Range("A88").Value = 45108
$ws.Range("A89").Value = 45139
$ws.Range("A90").Value = 45170
$ws.Range("A91").Value = 45200
$ws.Range("A92").Value = 45231
$ws.Range("A93").Value = 45261
$ws.Range("A94").Value = 45292
$ws.Range("A95").Value = 45323
$ws.Range("A96").Value = 45352
$ws.Range("A97").Value = 45383
$ws.Range("A98").Value = 45413
$ws.Range("A99").Value = 45444
$ws.Range("A100").Value = 45474
$ws.Range("A101").Value = 45505
$ws.Range("A102").Value = 45536
$ws.Range("A103").Value = 45566
$ws.Range("A104").Value = 45597
$ws.Range("A105").Value = 45627
$ws.Range("A106").Value = 45658
$ws.Range("A107").Value = 45689
$ws.Range("A108").Value = 45717
$ws.Range("A109").Value = 45748
$ws.Range("A110").Value = 45778
$ws.Range("A111").Value = 45809
$ws.Range("A112").Value = 45839
$ws.Range("A113").Value = 45870
$ws.Range("A114").Value = 45901
$ws.Range("A115").Value = 45931
$ws.Range("A116").Value = 45962
$ws.Range("A117").Value = 45992
$ws.Range("A118").Value = 46023
$ws.Range("A119").Value = 46054
$ws.Range("A120").Value = 46082
$ws.Range("A121").Value = 46113
$ws.Range("A122").Value = 46143
$ws.Range("A123").Value = 46174
$ws.Range("A124").Value = 46204
$ws.Range("A125").Value = 46235
$ws.Range("A126").Value = 46266
$ws.Range("A127").Value = 46296

# --- Footer: update Certified Correct By name/title ---
$ps = $ws.PageSetup
$ps.CenterFooter = '
CERTIFIED CORRECT BY: &UNANETTE B. SUSA&U
                                           OIC - HRMO'

# --- Selection: move active cell in bottom pane to K83 ---
$ws.Range("K83").Select()
